$wb = $excel.ActiveWorkbook

# --- Add the new "2024-2025" sheet -----------------------------------------
# Worksheets.Add() always inserts at the front with a default name ("Sheet1");
# rename + populate it BEFORE moving it into its final position, because after
# a Move() the original $ws reference resolves against a stale sheet index.
$ws = $wb.Worksheets.Add()
$ws.Name = "2024-2025"

# Column A width (approximates the ~36.71-char width used in the source file).
$ws.Columns("A").ColumnWidth = 35.9

# Header row (row 1) - same headers shared by every year tab.
$ws.Range("A1").Value = "Año fiscal/delitos"
$ws.Range("B1").Value = "Pendiente al inicio"
$ws.Range("C1").Value = "Casos presentados"
$ws.Range("D1").Value = "Casos a resolver"
$ws.Range("E1").Value = "Casos resueltos: Condenas"
$ws.Range("F1").Value = "Casos resueltos: Absoluciones"
$ws.Range("G1").Value = "Casos resueltos: Archivos"
$ws.Range("H1").Value = "Casos resueltos: Traslados"
$ws.Range("I1").Value = "Casos resueltos: Otros**"
$ws.Range("J1").Value = "Casos resueltos: Total"
$ws.Range("K1").Value = "Pendiente al final"

# Row 2 - fiscal year label for this tab.
$ws.Range("A2").Value = "2024-2025*"

# Row 3 - Acoso sexual
$ws.Range("A3").Value = "Acoso sexual"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 1
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 3

# Row 4 - Actos lascivos
$ws.Range("A4").Value = "Actos lascivos"
$ws.Range("B4").Value = 74
$ws.Range("C4").Value = 73
$ws.Range("D4").Value = 147
$ws.Range("E4").Value = 105
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 5
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 116
$ws.Range("K4").Value = 57

# Row 5 - Agresión sexual
$ws.Range("A5").Value = "Agresión sexual"
$ws.Range("B5").Value = 60
$ws.Range("C5").Value = 67
$ws.Range("D5").Value = 127
$ws.Range("E5").Value = 19
$ws.Range("G5").Value = 2
$ws.Range("J5").Value = 21
$ws.Range("K5").Value = 42

# Row 6 - Incesto
$ws.Range("A6").Value = "Incesto"
$ws.Range("B6").Value = 13
$ws.Range("C6").Value = 24
$ws.Range("D6").Value = 24
$ws.Range("E6").Value = 5
$ws.Range("J6").Value = 5
$ws.Range("K6").Value = 20

# Row 7 - Ley contra el acecho en Puerto Rico
$ws.Range("A7").Value = "Ley contra el acecho en Puerto Rico"
$ws.Range("B7").Value = 67
$ws.Range("C7").Value = 115
$ws.Range("D7").Value = 182
$ws.Range("E7").Value = 78
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = 21
$ws.Range("H7").Value = 6
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 110
$ws.Range("K7").Value = 76

# Row 8 - Tentativa de actos lascivos
$ws.Range("A8").Value = "Tentativa de actos lascivos"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 15
$ws.Range("J8").Value = 15
$ws.Range("K8").Value = 1

# Row 9 - Tentativa de agresión sexual
$ws.Range("A9").Value = "Tentativa de agresión sexual"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = 15
$ws.Range("G9").Value = 1
$ws.Range("J9").Value = 16

# Row 10 - Tentativa de incesto
$ws.Range("A10").Value = "Tentativa de incesto"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 2
$ws.Range("J10").Value = 2
$ws.Range("K10").Value = 1

# Row 11 - Tentativa de acoso sexual (new category, only present this year)
$ws.Range("A11").Value = "Tentativa de acoso sexual"
$ws.Range("E11").Value = 1
$ws.Range("J11").Value = 1

# Move the new tab to the end, right after "2023-2024".
$ws.Move($null, $wb.Worksheets.Item("2023-2024"))

# Re-fetch by name (Move() invalidates the old reference) and make it the
# active / selected sheet, matching activeTab="4" + the D6 selection.
$ws2 = $wb.Worksheets.Item("2024-2025")
$ws2.Activate()
$ws2.Range("D6").Select()
